$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "57.976.76"
$ws.Range("E2").Value = "  -0.24%  "

# Row 3
$ws.Range("D3").Value = "2.351.09"
$ws.Range("E3").Value = "  -0.42%  "

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("E4").Value = "  +0.29%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "541.25"
$ws.Range("E5").Value = "  -0.34%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "134.76"
$ws.Range("E6").Value = "  -0.38%  "

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.00"
$ws.Range("E7").Value = "  +0.25%  "

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.569"
$ws.Range("E8").Value = "  +6.04%  "

# Row 9
$ws.Range("E9").Value = "  +0.47%  "

# Row 10
$ws.Range("E10").Value = "  +2.67%  "

# Row 11
$ws.Range("E11").Value = "  -1.65%  "

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "23.84"
$ws.Range("E13").Value = "  +1.31%  "

# Row 14
$ws.Range("D14").Value = "2.769.78"
$ws.Range("E14").Value = "  +0.55%  "

# Row 15
$ws.Range("D15").Value = "57.901.92"
$ws.Range("E15").Value = "  -0.16%  "

# Row 16
$ws.Range("E16").Value = "  +0.82%  "

# Row 17
$ws.Range("D17").Value = "2.337.48"
$ws.Range("E17").Value = "  -0.74%  "

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "10.72"
$ws.Range("E18").Value = "  +1.40%  "

# Row 19
$ws.Range("E19").Value = "  +2.04%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "330.14"
$ws.Range("E20").Value = "  -2.50%  "

# Row 21
$ws.Range("E21").Value = "  -1.40%  "

# Row 22
$ws.Range("E22").Value = "  +0.47%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "62.80"
$ws.Range("E23").Value = "  +0.82%  "

# Row 24
$ws.Range("E24").Value = "  -2.48%  "

# Row 25
$ws.Range("E25").Value = "  +0.24%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "8.35"
$ws.Range("E26").Value = "  -1.78%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "1.34"
$ws.Range("E27").Value = "  -5.58%  "

# Row 28
$ws.Range("E28").Value = "  +0.12%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "170.00"
$ws.Range("E29").Value = "  -0.70%  "

# Row 30
$ws.Range("E30").Value = "  -0.18%  "

# Row 31
$ws.Range("E31").Value = "  -0.82%  "

# Row 32
$ws.Range("E32").Value = "  +1.16%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "18.38"
$ws.Range("E33").Value = "  -1.05%  "

# Row 34
$ws.Range("E34").Value = "  +0.09%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.00"
$ws.Range("E35").Value = "  +0.05%  "

# Row 36
$ws.Range("E36").Value = "  +1.19%  "

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.24"
$ws.Range("E37").Value = "  -1.63%  "

# Row 38
$ws.Range("E38").Value = "  -0.10%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "39.10"
$ws.Range("E39").Value = "  -0.67%  "

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "142.88"
$ws.Range("E40").Value = "  -3.95%  "

# Row 41
$ws.Range("B41").Value = "Filecoin"
$ws.Range("C41").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "3.65"
$ws.Range("E41").Value = "  +0.59%  "

# Row 42
$ws.Range("B42").Value = "PolygonEcosystemToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.378"
$ws.Range("E42").Value = "  -0.22%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "288.92"
$ws.Range("E43").Value = "  +1.27%  "

# Row 45
$ws.Range("E45").Value = "  +0.52%  "

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "19.14"
$ws.Range("E46").Value = "  -0.60%  "

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.567"
$ws.Range("E47").Value = "  +1.48%  "

# Row 48
$ws.Range("E48").Value = "  +1.52%  "

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.385"
$ws.Range("E49").Value = "  +0.56%  "

# Row 50
$ws.Range("E50").Value = "  +0.58%  "

# Row 51
$ws.Range("E51").Value = "  -0.73%  "
